$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "login--functionality;login-with-valid-username-and-password"
$ws.Range("B12").Value = "failed"
$ws.Range("C12").Value = "chrome"
$ws.Range("D12").Value = "16.09.21"
